$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency price/volume table with the latest scraped values
# (scheduled GitHub Actions refresh). Rows 45/46 (EnergySwap/Decentraland)
# swapped ranking order along with their data.
# D-column prices are momentarily forced to text format before assignment so
# that numeric-looking strings (e.g. "1.001", "13.10") keep their exact
# textual representation instead of being auto-converted to numbers, then the
# cell style is restored to Normal so no visible formatting change remains.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.923.91'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.867.02'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.46%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.99%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4974'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3803'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08941'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -9.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.117'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.44'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.93%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.298'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.61'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.858.21'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.214'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001097'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06620'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.86'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.074'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.942.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.37'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.284'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.67%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.075.59'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.518'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.81'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.65'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.85%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1053'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.053'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.574'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.590'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.326'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06529'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02403'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2180'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.274'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.198'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.63'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6356'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.889'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.05%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5984'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.90%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.10'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.282'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.74%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.667'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.217'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.962'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.72'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.21%  '
